{"js": "// Update the date line and the multiplication problems to the new values.\nconst replacements = [\n  [\"2026-01-15 Thursday\", \"2026-01-16 Friday\"],\n  [\"743\u00d78=\", \"426\u00d75=\"],\n  [\"242\u00d72=\", \"541\u00d74=\"],\n  [\"591\u00d72=\", \"575\u00d77=\"],\n  [\"344\u00d73=\", \"832\u00d79=\"],\n  [\"607\u00d76=\", \"924\u00d72=\"],\n  [\"183\u00d76=\", \"211\u00d77=\"],\n  [\"629\u00d72=\", \"843\u00d76=\"],\n  [\"610\u00d75=\", \"992\u00d72=\"],\n  [\"889\u00d77=\", \"880\u00d78=\"],\n  [\"976\u00d79=\", \"456\u00d79=\"],\n  [\"625\u00d77=\", \"631\u00d77=\"],\n  [\"365\u00d75=\", \"820\u00d74=\"],\n  [\"956\u00d78=\", \"239\u00d76=\"],\n  [\"220\u00d73=\", \"502\u00d76=\"],\n  [\"119\u00d79=\", \"762\u00d79=\"],\n  [\"819\u00d78=\", \"416\u00d72=\"],\n  [\"766\u00d73=\", \"221\u00d76=\"],\n  [\"466\u00d77=\", \"533\u00d79=\"],\n  [\"602\u00d72=\", \"421\u00d74=\"],\n  [\"125\u00d73=\", \"545\u00d78=\"],\n  [\"673\u00d73=\", \"142\u00d74=\"],\n  [\"277\u00d77=\", \"304\u00d77=\"],\n  [\"285\u00d79=\", \"364\u00d75=\"],\n  [\"151\u00d77=\", \"864\u00d76=\"],\n  [\"245\u00d75=\", \"645\u00d77=\"],\n];\n\nfor (const [oldText, newText] of replacements) {\n  const results = context.document.body.search(oldText, { matchCase: true });\n  results.load(\"items\");\n  await context.sync();\n\n  for (const item of results.items) {\n    item.insertText(newText, Word.InsertLocation.replace);\n  }\n  await context.sync();\n}\n", "ps1": "$d = $word.ActiveDocument\n\n$replacements = @(\n    @(\"2026-01-15 Thursday\", \"2026-01-16 Friday\"),\n    @(\"743\u00d78=\", \"426\u00d75=\"),\n    @(\"242\u00d72=\", \"541\u00d74=\"),\n    @(\"591\u00d72=\", \"575\u00d77=\"),\n    @(\"344\u00d73=\", \"832\u00d79=\"),\n    @(\"607\u00d76=\", \"924\u00d72=\"),\n    @(\"183\u00d76=\", \"211\u00d77=\"),\n    @(\"629\u00d72=\", \"843\u00d76=\"),\n    @(\"610\u00d75=\", \"992\u00d72=\"),\n    @(\"889\u00d77=\", \"880\u00d78=\"),\n    @(\"976\u00d79=\", \"456\u00d79=\"),\n    @(\"625\u00d77=\", \"631\u00d77=\"),\n    @(\"365\u00d75=\", \"820\u00d74=\"),\n    @(\"956\u00d78=\", \"239\u00d76=\"),\n    @(\"220\u00d73=\", \"502\u00d76=\"),\n    @(\"119\u00d79=\", \"762\u00d79=\"),\n    @(\"819\u00d78=\", \"416\u00d72=\"),\n    @(\"766\u00d73=\", \"221\u00d76=\"),\n    @(\"466\u00d77=\", \"533\u00d79=\"),\n    @(\"602\u00d72=\", \"421\u00d74=\"),\n    @(\"125\u00d73=\", \"545\u00d78=\"),\n    @(\"673\u00d73=\", \"142\u00d74=\"),\n    @(\"277\u00d77=\", \"304\u00d77=\"),\n    @(\"285\u00d79=\", \"364\u00d75=\"),\n    @(\"151\u00d77=\", \"864\u00d76=\"),\n    @(\"245\u00d75=\", \"645\u00d77=\")\n)\n\nforeach ($pair in $replacements) {\n    $range = $d.Content\n    $range.Find.ClearFormatting()\n    $range.Find.Replacement.ClearFormatting()\n    $range.Find.Text = $pair[0]\n    $range.Find.Replacement.Text = $pair[1]\n    $range.Find.Execute($pair[0], $false, $false, $false, $false, $false, $true, 1, $false, $pair[1], 2)\n}\n"}
